# Split the single "Coming soon" continuation run into three runs:
#   " How to "  |  "digitally "  |  "implement a biquadratic notch filter"
# (was one run: " How to implement a biquadratic notch filter digitally / in code")

$d = $word.ActiveDocument

# Locate the existing text precisely (robust to any offset assumptions).
$target = $d.Content
$target.Find.Execute(" How to implement a biquadratic notch filter digitally / in code", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Rewrite the found range in place with the first chunk's text. This keeps
# the run "alive" (non-empty) so it is not merged away into its bold
# neighbor, and gives us a known end position to build on.
$target.Text = " How to "

# Insert the second chunk right after, as its own run.
$pos2 = $target.End
$ins2 = $d.Range($pos2, $pos2)
$ins2.InsertAfter("digitally ")

# Insert the third chunk right after that, as its own run.
$pos3 = $pos2 + 10
$ins3 = $d.Range($pos3, $pos3)
$ins3.InsertAfter("implement a biquadratic notch filter")
